# Topic 2 lease example - add "% of MLP" column (K) and "Income before taxes" row (17)
# plus a new "Steps:" explanatory text box, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column K: "% of MLP" = Total Expense (J) / MLP (C), for each year row
# ---------------------------------------------------------------------------

# Header cell K3, styled like the other header cells in row 3 (bold-ish box w/ border)
$ws.Range("K3").Value = "% of MLP"
$ws.Range("K3").HorizontalAlignment = -4108   # xlCenter
$ws.Range("K3").VerticalAlignment = -4160     # xlTop
$ws.Range("K3").Borders.Item(9).LineStyle = 1 # xlEdgeBottom thin border
$ws.Range("K3").Font.Bold = $false

# Blank formatted cells above/below the header, matching neighbouring columns
$ws.Range("K1").Borders.Item(9).LineStyle = 1
$ws.Range("K2").Borders.Item(9).LineStyle = 1
$ws.Range("K4").HorizontalAlignment = -4108

# Data rows 5-12: K = J/C, percentage, no border
for ($r = 5; $r -le 12; $r++) {
    $cell = $ws.Range("K$r")
    $cell.Formula = "=J$r/C$r"
    $cell.NumberFormat = "0%"
    $cell.HorizontalAlignment = -4108
}

# Row 13 is the last detail row and carries the bottom border of the table
$ws.Range("K13").Formula = "=J13/C13"
$ws.Range("K13").NumberFormat = "0%"
$ws.Range("K13").HorizontalAlignment = -4108
$ws.Range("K13").Borders.Item(9).LineStyle = 1

# Row 14 (Total row) - leave K blank but keep it inside the shaded/bordered block
$ws.Range("K14").Interior.ColorIndex = 6
$ws.Range("K14").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# 2. New row 17: "Income before taxes" ratio, inserted above the existing
#    "Total asset" / "Total liability" rows
# ---------------------------------------------------------------------------

$ws.Range("B17").Value = "Income before taxes"
$ws.Range("C17").Value = 2660
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("E17").Formula = "=(J9-C9)/C17"
$ws.Range("E17").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 3. Explanatory "Steps:" text box describing the workflow
# ---------------------------------------------------------------------------

$shp = $ws.Shapes.AddTextbox(1, 63.5, 312.5, 473.5, 137.5)
$shp.Name = "Text Frame 1"

$steps = "Steps: " + [char]10 +
  "1. What is the implicit interest rate of the operating lease given the PV?" + [char]10 +
  "2. Converting MLP into a liability with two components: (1) interest + (2) capital amortization (as it would be a normal credit liability). Here we keep the lease jergoen an call it " + [char]8220 + "Lease Obligation" + [char]8221 + "." + [char]10 +
  "3. Estimate the value of the operating lease asset." + [char]10 +
  "For simplicity, we assume that its value equals the PV of the MLP (liability). With this we can compute depreciation expense of the respective assets." + [char]10 +
  "4. Income statement: Compute the interest expense and depreciation expense for each year." + [char]10 +
  ""

$shp.TextFrame.Characters().Text = $steps
$shp.TextFrame.Characters().Font.Name = "Times New Roman"
$shp.TextFrame.Characters().Font.Size = 12
$shp.TextFrame.Characters().Font.Bold = $false

# ---------------------------------------------------------------------------
# 4. Selection cosmetic change recorded in the saved file
# ---------------------------------------------------------------------------
$ws.Range("O13").Select()
